$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text format
# first, so Excel stores them as strings (matching the original inlineStr layout)
# instead of silently converting them to numeric cells.
$ws.Range("D2").Value = "97.214.67"
$ws.Range("E2").Value = "  +5.09%  "
$ws.Range("D3").Value = "3.116.62"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.66"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "609.82"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("E7").Value = "  +1.98%  "
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "3.114.79"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.783"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "96.656.16"
$ws.Range("E13").Value = "  +4.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000240"
$ws.Range("E14").Value = "  -1.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.70"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.41"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").Value = "3.691.08"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "3.122.10"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "530.73"
$ws.Range("E19").Value = "  +20.99%  "
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.53"
$ws.Range("E20").Value = "  -6.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.48"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.65"
$ws.Range("E22").Value = "  -3.29%  "
$ws.Range("E23").Value = "  -4.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.80"
$ws.Range("E24").Value = "  -3.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.46"
$ws.Range("E25").Value = "  -1.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.27"
$ws.Range("E26").Value = "  +3.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.54"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "3.280.94"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.235"
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("E31").Value = "  -4.25%  "
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.98"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.74"
$ws.Range("E34").Value = "  +4.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.867"
$ws.Range("E35").Value = "  -13.25%  "
$ws.Range("E36").Value = "  -8.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.23"
$ws.Range("E37").Value = "  -10.72%  "
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.21"
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "477.91"
$ws.Range("E40").Value = "  +3.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.436"
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("E42").Value = "  -3.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.59"
$ws.Range("E43").Value = "  -9.88%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.16"
$ws.Range("E45").Value = "  -4.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "160.88"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.90"
$ws.Range("E47").Value = "  +3.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.689"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.48"
$ws.Range("E49").Value = "  +3.30%  "
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.43"
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("E51").Value = "  +0.01%  "
